# This script applies numeric updates to the "two-digit number divided by
# one-digit number" worksheet table. Each division problem's dividend/divisor
# pair is replaced with a new pair, per the target revision.
#
# Replacements are executed in document order (top-to-bottom, left-to-right)
# so that a freshly written "new" value that happens to equal some other
# cell's "old" value is never re-matched by a later Find/Replace call
# (Word's Find.Execute with wdReplaceAll rescans the whole $d.Content range
# on every call, not just from the cursor position).

$d = $word.ActiveDocument

$d.Content.Find.Execute("92÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "61÷3=", 2) | Out-Null
$d.Content.Find.Execute("98÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "31÷3=", 2) | Out-Null
$d.Content.Find.Execute("61÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "75÷6=", 2) | Out-Null
$d.Content.Find.Execute("88÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "15÷6=", 2) | Out-Null
$d.Content.Find.Execute("12÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "98÷9=", 2) | Out-Null
$d.Content.Find.Execute("22÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷5=", 2) | Out-Null
$d.Content.Find.Execute("41÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "87÷5=", 2) | Out-Null
$d.Content.Find.Execute("71÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "20÷5=", 2) | Out-Null
$d.Content.Find.Execute("61÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "61÷9=", 2) | Out-Null
$d.Content.Find.Execute("79÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "26÷5=", 2) | Out-Null
$d.Content.Find.Execute("40÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "51÷3=", 2) | Out-Null
$d.Content.Find.Execute("65÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "37÷4=", 2) | Out-Null
$d.Content.Find.Execute("73÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "55÷5=", 2) | Out-Null
$d.Content.Find.Execute("87÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "42÷4=", 2) | Out-Null
$d.Content.Find.Execute("74÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "59÷8=", 2) | Out-Null
$d.Content.Find.Execute("43÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "51÷9=", 2) | Out-Null
$d.Content.Find.Execute("28÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "31÷3=", 2) | Out-Null
$d.Content.Find.Execute("35÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "72÷3=", 2) | Out-Null
$d.Content.Find.Execute("50÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "81÷3=", 2) | Out-Null
$d.Content.Find.Execute("24÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "67÷6=", 2) | Out-Null
$d.Content.Find.Execute("21÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "66÷3=", 2) | Out-Null
$d.Content.Find.Execute("33÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "71÷9=", 2) | Out-Null
$d.Content.Find.Execute("99÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "14÷6=", 2) | Out-Null
$d.Content.Find.Execute("23÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "41÷6=", 2) | Out-Null
$d.Content.Find.Execute("66÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "47÷2=", 2) | Out-Null
